$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels:
#  B1: "TIPO DE PROMO"     -> "TIPO_PROMO"
#  D1: "APPLICATION FORM"  -> "APPLICATION_FORM"
# (other header cells / A2 stay the same, just re-indexed internally)
$ws.Range("B1").Value = "TIPO_PROMO"
$ws.Range("D1").Value = "APPLICATION_FORM"
